$d = $word.ActiveDocument

# The document currently has a single paragraph:
#   "ARQUIVO TESTELM V.1"  (followed by the hidden _GoBack bookmark)
# We need to split it into two paragraphs, keeping the first paragraph's
# text as-is, and adding a new second paragraph with the bold / 20pt /
# underlined title text. The _GoBack bookmark must end up at the end of
# the new second paragraph (immediately after its text, before its own
# paragraph mark) - i.e. where Word would naturally leave it after typing.

$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphAfter()

# This is now the new (second, currently empty) paragraph.
$newPara = $d.Paragraphs.Item(2)

# Type the new title text. A temporary trailing marker character ("X") is
# appended so that, a moment from now, we can anchor a zero-length
# bookmark right after the real text without landing exactly on a
# paragraph-mark boundary (collapsed bookmarks placed precisely at a
# paragraph mark / content end are not positioned reliably).
$newPara.Range.Text = "“Meu primeiro exercício com GIT Desktop”X"

# Apply the run/paragraph formatting: bold, 20pt (sz is in half-points,
# so 40), single underline.
$newPara.Range.Font.Bold = $true
$newPara.Range.Font.Size = 20
$newPara.Range.Font.Underline = 1

# Re-fetch the paragraph range after the edits above.
$newPara = $d.Paragraphs.Item(2)
$endPos = $newPara.Range.End - 2   # position right before the temp "X"

# Relocate the _GoBack bookmark to sit right after the new text.
$d.Bookmarks.Item("_GoBack").Delete()
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the temporary marker character now that the bookmark is anchored.
$markerRange = $d.Range($endPos, $endPos + 1)
$markerRange.Delete()
